$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (so the new "CRUCETA" entry becomes row 4,
# shifting the old row4 "CRUZ PAMPA YAPATERA" down to row 5, etc.)
$ws.Rows.Item(5).Insert()

# Update cells for rows 2-15 with final data
$data = @(
    @(2, "PIURA", "19 DE AGOSTO", "PIURA", 5840, 370, 6.335616438356165),
    @(3, "LA LIBERTAD", "CIUDAD DE DIOS", "PACASMAYO", 3000, 248, 8.266666666666666),
    @(4, "PIURA", "CRUCETA", "PIURA", 4410, 176, 3.990929705215419),
    @(5, "PIURA", "CRUZ PAMPA YAPATERA", "MORROPON", 5000, 229, 4.58),
    @(6, "LA LIBERTAD", "CURVA DE SUN", "TRUJILLO", 8000, 191, 2.3875),
    @(7, "CAJAMARCA", "HUAMBOCANCHA ALTA", "CAJAMARCA", 5894, 287, 4.86935866983373),
    @(8, "LA LIBERTAD", "HUANCHAQUITO", "TRUJILLO", 7000, 208, 2.971428571428572),
    @(9, "LAMBAYEQUE", "LA COLORADA", "LAMBAYEQUE", 4000, 429, 10.725),
    @(10, "PIURA", "LA PEÑITA", "PIURA", 5868, 246, 4.192229038854806),
    @(11, "PIURA", "MALINGAS", "PIURA", 5574, 350, 6.279153211338357),
    @(12, "CAJAMARCA", "OTUZCO", "CAJAMARCA", 6000, 149, 2.483333333333333),
    @(13, "CAJAMARCA", "SAN ANTONIO BAJO", "HUALGAYOC", 5000, 181, 3.62),
    @(14, "PIURA", "LA VILLA LETIRA - BECARA", "SECHURA", 6142, 860, 14.0019537609899),
    @(15, "PIURA", "VIVIATE", "PAITA", 5025, 537, 10.6865671641791)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
